$wb = $excel.ActiveWorkbook

# --- "provenance" sheet updates ---
$ws = $wb.Worksheets.Item("provenance")

# python source file: tools_xl.pyc -> tools_xl.py
$ws.Range("B4").Value = "tools_xl.py"

# directory: /Users/dantopa/... -> /Volumes/Tlaltecuhtli/...
$ws.Range("B5").Value = "/Volumes/Tlaltecuhtli/repos/GitHub/topa-development/amanzi/aqua/xl"

# python version string (Python 2.7 -> Python 3.7, new build string)
$ws.Range("B6").Value = "3.7.0 (default, Jun 28 2018, 07:39:16) `n[Clang 4.0.1 (tags/RELEASE_401/final)]"

# $USER: dantopa -> l127914
$ws.Range("B9").Value = "l127914"

# $HOSTNAME: MacBookPro11,3 -> Cauchy.Schwarz
$ws.Range("B10").Value = "Cauchy.Schwarz"

# $HOME: /Users/dantopa -> /Users/l127914
$ws.Range("B11").Value = "/Users/l127914"

# timestamp value (serial date/time number)
$ws.Range("B12").Value = 43437.44765421725

# --- "08-BC" sheet updates ---
$ws2 = $wb.Worksheets.Item("08-BC")

$ws2.Range("B3").Value = "1.assigned_region"
$ws2.Range("B4").Value = "2.liquid_phase"
$ws2.Range("B5").Value = "3.solid_phase"
